$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from row 15 to row 16 first, to match styles (date/time/text format + borders)
$ws.Range("A15:C15").Copy()
$ws.Range("A16:C16").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Add new row 16 data (use raw numeric serials so Excel doesn't auto-apply a date/time number format)
$ws.Range("A16").Value2 = 43916
$ws.Range("B16").Value2 = 0.33333333333333331
$ws.Range("C16").Value = "Initiation aux tests fonctionnels"

# Update selection to match diff (activeCell B17)
$ws.Range("B17").Select()
